$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33 (diff @ line 2249)
$ws.Cells.Item(33, 8).Value = 875.96
$ws.Cells.Item(33, 9).Value = 872.8946999999999
$ws.Cells.Item(33, 10).Value = 885.6667
$ws.Cells.Item(33, 11).Value = 872.8946999999999
$ws.Cells.Item(33, 12).Value = 885.6667
$ws.Cells.Item(33, 13).Value = -643.8946999999999
$ws.Cells.Item(33, 14).Value = -1343.6667
# Row 70 (diff @ line 4080)
$ws.Cells.Item(70, 8).Value = 12894.571
$ws.Cells.Item(70, 10).Value = 12894.571
$ws.Cells.Item(70, 12).Value = 38683.713
$ws.Cells.Item(70, 14).Value = -39223.713
# Row 73 (diff @ line 4227)
$ws.Cells.Item(73, 8).Value = 12894.571
$ws.Cells.Item(73, 10).Value = 12894.571
$ws.Cells.Item(73, 12).Value = 38683.713
$ws.Cells.Item(73, 14).Value = -40555.713
# Row 86 (diff @ line 4876)
$ws.Cells.Item(86, 8).Value = 2147.4614
$ws.Cells.Item(86, 9).Value = 2000
$ws.Cells.Item(86, 11).Value = 2000
$ws.Cells.Item(86, 13).Value = -877
# Row 89 (diff @ line 5032)
$ws.Cells.Item(89, 8).Value = 2147.4614
$ws.Cells.Item(89, 9).Value = 2000
$ws.Cells.Item(89, 11).Value = 10000
$ws.Cells.Item(89, 13).Value = -4384
# Row 97 (diff @ line 5436)
$ws.Cells.Item(97, 8).Value = 784.36365
$ws.Cells.Item(97, 10).Value = 777
$ws.Cells.Item(97, 12).Value = 2331
$ws.Cells.Item(97, 14).Value = -3323
# Row 112 (diff @ line 6195)
$ws.Cells.Item(112, 8).Value = 5017.2285
$ws.Cells.Item(112, 10).Value = 5237.4243
$ws.Cells.Item(112, 12).Value = 15712.2729
$ws.Cells.Item(112, 14).Value = -17928.2729
# Row 116 (diff @ line 6397)
$ws.Cells.Item(116, 8).Value = 5252.625
$ws.Cells.Item(116, 9).Value = 4701.7144
$ws.Cells.Item(116, 11).Value = 4701.7144
$ws.Cells.Item(116, 13).Value = -1259.7144
# Row 138 (diff @ line 7487)
$ws.Cells.Item(138, 8).Value = 2684.1082
$ws.Cells.Item(138, 9).Value = 1370.2106
$ws.Cells.Item(138, 10).Value = 4071
$ws.Cells.Item(138, 11).Value = 4110.6318
$ws.Cells.Item(138, 12).Value = 12213
$ws.Cells.Item(138, 13).Value = 1029.3682
$ws.Cells.Item(138, 14).Value = -22493
# Row 141 (diff @ line 7637)
$ws.Cells.Item(141, 8).Value = 1000
$ws.Cells.Item(141, 9).Value = 1000
$ws.Cells.Item(141, 11).Value = 3000
$ws.Cells.Item(141, 13).Value = 2180

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 97 (diff @ line 12414)
$ws.Cells.Item(97, 8).Value = 1624.75
$ws.Cells.Item(97, 9).Value = 1104.125
$ws.Cells.Item(97, 11).Value = 1104.125
$ws.Cells.Item(97, 13).Value = -608.125
# Row 102 (diff @ line 12653)
$ws.Cells.Item(102, 8).Value = 4554.9
$ws.Cells.Item(102, 9).Value = 4005.158
$ws.Cells.Item(102, 11).Value = 4005.158
$ws.Cells.Item(102, 13).Value = -2383.158
# Row 122 (diff @ line 13618)
$ws.Cells.Item(122, 8).Value = 3999
$ws.Cells.Item(122, 9).Value = 3998.5
$ws.Cells.Item(122, 11).Value = 11995.5
$ws.Cells.Item(122, 13).Value = -9545.5
# Row 132 (diff @ line 14111)
$ws.Cells.Item(132, 8).Value = 3787.6453
$ws.Cells.Item(132, 9).Value = 1284.2
$ws.Cells.Item(132, 10).Value = 4269.077
$ws.Cells.Item(132, 11).Value = 3852.6
$ws.Cells.Item(132, 12).Value = 12807.231
$ws.Cells.Item(132, 13).Value = -1322.6
$ws.Cells.Item(132, 14).Value = -17867.231

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94 (diff @ line 19197)
$ws.Cells.Item(94, 8).Value = 946.8
$ws.Cells.Item(94, 9).Value = 958.5
$ws.Cells.Item(94, 10).Value = 900
$ws.Cells.Item(94, 11).Value = 958.5
$ws.Cells.Item(94, 12).Value = 900
$ws.Cells.Item(94, 13).Value = -507.5
$ws.Cells.Item(94, 14).Value = -1802
# Row 99 (diff @ line 19442)
$ws.Cells.Item(99, 8).Value = 2721.7778
$ws.Cells.Item(99, 9).Value = 1749.8334
$ws.Cells.Item(99, 11).Value = 1749.8334
$ws.Cells.Item(99, 13).Value = -251.8334
# Row 107 (diff @ line 19837)
$ws.Cells.Item(107, 8).Value = 31023.354
$ws.Cells.Item(107, 9).Value = 40069.08
$ws.Cells.Item(107, 11).Value = 40069.08
$ws.Cells.Item(107, 13).Value = -38149.08
# Row 122 (diff @ line 20551)
$ws.Cells.Item(122, 8).Value = 50390
$ws.Cells.Item(122, 10).Value = 50390
$ws.Cells.Item(122, 12).Value = 50390
$ws.Cells.Item(122, 14).Value = -60190
# Row 134 (diff @ line 21130)
$ws.Cells.Item(134, 8).Value = 2302.8096
$ws.Cells.Item(134, 9).Value = 2050.5881
$ws.Cells.Item(134, 11).Value = 6151.7643
$ws.Cells.Item(134, 13).Value = -3616.7643

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7 (diff @ line 21867)
$ws.Cells.Item(7, 8).Value = 59
$ws.Cells.Item(7, 9).Value = 59
$ws.Cells.Item(7, 11).Value = 59
$ws.Cells.Item(7, 13).Value = 54
# Row 16 (diff @ line 22305)
$ws.Cells.Item(16, 8).Value = 997.4286
$ws.Cells.Item(16, 9).Value = 843.61536
$ws.Cells.Item(16, 11).Value = 843.61536
$ws.Cells.Item(16, 13).Value = -556.61536
# Row 31 (diff @ line 23052)
$ws.Cells.Item(31, 8).Value = 4763572.5
$ws.Cells.Item(31, 9).Value = 5557027.5
$ws.Cells.Item(31, 11).Value = 5557027.5
$ws.Cells.Item(31, 13).Value = -5556732.5
# Row 34 (diff @ line 23199)
$ws.Cells.Item(34, 8).Value = 4763572.5
$ws.Cells.Item(34, 9).Value = 5557027.5
$ws.Cells.Item(34, 11).Value = 5557027.5
$ws.Cells.Item(34, 13).Value = -5556825.5
# Row 62 (diff @ line 24556)
$ws.Cells.Item(62, 8).Value = 4798.2
$ws.Cells.Item(62, 10).Value = 4798.2
$ws.Cells.Item(62, 12).Value = 4798.2
$ws.Cells.Item(62, 14).Value = -6046.2
# Row 65 (diff @ line 24703)
$ws.Cells.Item(65, 8).Value = 4798.2
$ws.Cells.Item(65, 10).Value = 4798.2
$ws.Cells.Item(65, 12).Value = 23991
$ws.Cells.Item(65, 14).Value = -30231
# Row 94 (diff @ line 26139)
$ws.Cells.Item(94, 8).Value = 1800
$ws.Cells.Item(94, 10).Value = 1800
$ws.Cells.Item(94, 12).Value = 1800
$ws.Cells.Item(94, 14).Value = -2702
# Row 99 (diff @ line 26384)
$ws.Cells.Item(99, 8).Value = 11911.637
$ws.Cells.Item(99, 9).Value = 10198
$ws.Cells.Item(99, 11).Value = 10198
$ws.Cells.Item(99, 13).Value = -8700
# Row 113 (diff @ line 27070)
$ws.Cells.Item(113, 8).Value = 997.4286
$ws.Cells.Item(113, 9).Value = 843.61536
$ws.Cells.Item(113, 11).Value = 843.61536
$ws.Cells.Item(113, 13).Value = 1326.38464
# Row 126 (diff @ line 27707)
$ws.Cells.Item(126, 8).Value = 11911.637
$ws.Cells.Item(126, 9).Value = 10198
$ws.Cells.Item(126, 11).Value = 30594
$ws.Cells.Item(126, 13).Value = -28124
# Row 134 (diff @ line 28108)
$ws.Cells.Item(134, 8).Value = 2617.7917
$ws.Cells.Item(134, 9).Value = 2215.3684
$ws.Cells.Item(134, 11).Value = 6646.1052
$ws.Cells.Item(134, 13).Value = -4111.1052

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 6 (diff @ line 28796)
$ws.Cells.Item(6, 8).Value = 5564.1665
$ws.Cells.Item(6, 9).Value = 4678.2
$ws.Cells.Item(6, 10).Value = 9994
$ws.Cells.Item(6, 11).Value = 14034.6
$ws.Cells.Item(6, 12).Value = 29982
$ws.Cells.Item(6, 13).Value = -13921.6
$ws.Cells.Item(6, 14).Value = -30208
# Row 11 (diff @ line 29044)
$ws.Cells.Item(11, 8).Value = 202245.8
$ws.Cells.Item(11, 10).Value = 10000
$ws.Cells.Item(11, 12).Value = 30000
$ws.Cells.Item(11, 14).Value = -30280
# Row 26 (diff @ line 29803)
$ws.Cells.Item(26, 8).Value = 5011
$ws.Cells.Item(26, 9).Value = 20
$ws.Cells.Item(26, 10).Value = 10002
$ws.Cells.Item(26, 11).Value = 60
$ws.Cells.Item(26, 12).Value = 30006
$ws.Cells.Item(26, 13).Value = 228
$ws.Cells.Item(26, 14).Value = -30582
# Row 52 (diff @ line 31092)
$ws.Cells.Item(52, 8).Value = 2015.5
$ws.Cells.Item(52, 10).Value = 2015.5
$ws.Cells.Item(52, 12).Value = 6046.5
$ws.Cells.Item(52, 14).Value = -6578.5
# Row 56 (diff @ line 31291)
$ws.Cells.Item(56, 8).Value = 7775.4
$ws.Cells.Item(56, 9).Value = 7775.4
$ws.Cells.Item(56, 11).Value = 7775.4
$ws.Cells.Item(56, 13).Value = -7245.4
# Row 68 (diff @ line 31894)
$ws.Cells.Item(68, 8).Value = 4274.129
$ws.Cells.Item(68, 9).Value = 827.6667
$ws.Cells.Item(68, 10).Value = 4643.393
$ws.Cells.Item(68, 11).Value = 2483.0001
$ws.Cells.Item(68, 12).Value = 13930.179
$ws.Cells.Item(68, 13).Value = -1672.0001
$ws.Cells.Item(68, 14).Value = -15552.179
# Row 71 (diff @ line 32050)
$ws.Cells.Item(71, 8).Value = 4274.129
$ws.Cells.Item(71, 9).Value = 827.6667
$ws.Cells.Item(71, 10).Value = 4643.393
$ws.Cells.Item(71, 11).Value = 7449.0003
$ws.Cells.Item(71, 12).Value = 41790.537
$ws.Cells.Item(71, 13).Value = -3393.0003
$ws.Cells.Item(71, 14).Value = -49902.537
# Row 132 (diff @ line 35153)
$ws.Cells.Item(132, 8).Value = 1369.7059
$ws.Cells.Item(132, 9).Value = 1688.3
$ws.Cells.Item(132, 11).Value = 15194.7
$ws.Cells.Item(132, 13).Value = -12664.7

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80 (diff @ line 39556)
$ws.Cells.Item(80, 8).Value = 11928.571
$ws.Cells.Item(80, 10).Value = 30000
$ws.Cells.Item(80, 12).Value = 30000
$ws.Cells.Item(80, 14).Value = -31996
# Row 83 (diff @ line 39703)
$ws.Cells.Item(83, 8).Value = 11928.571
$ws.Cells.Item(83, 10).Value = 30000
$ws.Cells.Item(83, 12).Value = 150000
$ws.Cells.Item(83, 14).Value = -159984
# Row 126 (diff @ line 41783)
$ws.Cells.Item(126, 8).Value = 3334.5
$ws.Cells.Item(126, 10).Value = 3817.4
$ws.Cells.Item(126, 12).Value = 11452.2
$ws.Cells.Item(126, 14).Value = -16392.2
# Row 132 (diff @ line 42080)
$ws.Cells.Item(132, 8).Value = 4567.4546
$ws.Cells.Item(132, 9).Value = 4088
$ws.Cells.Item(132, 11).Value = 12264
$ws.Cells.Item(132, 13).Value = -9734

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61 (diff @ line 45591)
$ws.Cells.Item(61, 8).Value = 1033.7727
$ws.Cells.Item(61, 9).Value = 795.65
$ws.Cells.Item(61, 11).Value = 795.65
$ws.Cells.Item(61, 13).Value = -593.65
# Row 104 (diff @ line 47716)
$ws.Cells.Item(104, 8).Value = 10286.25
$ws.Cells.Item(104, 10).Value = 10286.25
$ws.Cells.Item(104, 12).Value = 10286.25
$ws.Cells.Item(104, 14).Value = -17274.25
# Row 113 (diff @ line 48154)
$ws.Cells.Item(113, 8).Value = 1033.7727
$ws.Cells.Item(113, 9).Value = 795.65
$ws.Cells.Item(113, 11).Value = 795.65
$ws.Cells.Item(113, 13).Value = 1374.35
# Row 132 (diff @ line 49085)
$ws.Cells.Item(132, 8).Value = 4720.143
$ws.Cells.Item(132, 9).Value = 4720.143
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 14160.429
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -11630.429
$ws.Cells.Item(132, 14).ClearContents()
# Row 136 (diff @ line 49281)
$ws.Cells.Item(136, 8).Value = 6399.2915
$ws.Cells.Item(136, 9).Value = 6708.625
$ws.Cells.Item(136, 11).Value = 20125.875
$ws.Cells.Item(136, 13).Value = -17575.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107 (diff @ line 54772)
$ws.Cells.Item(107, 8).Value = 1381.625
$ws.Cells.Item(107, 9).Value = 1544.1666
$ws.Cells.Item(107, 10).Value = 894
$ws.Cells.Item(107, 11).Value = 4632.4998
$ws.Cells.Item(107, 12).Value = 2682
$ws.Cells.Item(107, 13).Value = -2712.4998
$ws.Cells.Item(107, 14).Value = -6522
# Row 122 (diff @ line 55498)
$ws.Cells.Item(122, 8).Value = 180283.08
$ws.Cells.Item(122, 9).Value = 211061.64
$ws.Cells.Item(122, 11).Value = 633184.92
$ws.Cells.Item(122, 13).Value = -630734.92
# Row 132 (diff @ line 55988)
$ws.Cells.Item(132, 8).Value = 3281.4443
$ws.Cells.Item(132, 9).Value = 4189.2354
$ws.Cells.Item(132, 11).Value = 12567.7062
$ws.Cells.Item(132, 13).Value = -10037.7062

Write-Host "Applied Midgardsormr_Profits updates across all sheets."